$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation date (2021-02-08, serial 44235) was missing from the
# series and has now been added between the existing rows for 44234 and
# 44236. Insert a fresh row at 93 so everything below shifts down by one
# (dates/format carry down automatically with the insert).
$ws.Rows("93:93").Insert()

# The freshly inserted row has no style yet - copy the date cell format
# from the row above (A92) so A93 keeps the same date number format.
$ws.Range("A92").Copy()
$ws.Range("A93").PasteSpecial(-4122)

# Fill in the newly inserted row for 2021-02-08.
$ws.Range("A93").Value = 44235
$ws.Range("B93").Value = 8
$ws.Range("C93").Value = 32
$ws.Range("D93").Value = 283.4617769510143

# The rolling-window columns (C: 7-day sum, D: sum per 100k inhabitants)
# for the days surrounding the newly inserted date need to be
# recalculated to include it.
$ws.Range("C90").Value = 42
$ws.Range("D90").Value = 372.0435822482062

$ws.Range("C91").Value = 35
$ws.Range("D91").Value = 310.0363185401719

$ws.Range("C92").Value = 30
$ws.Range("D92").Value = 265.7454158915759

$ws.Range("C94").Value = 35
$ws.Range("D94").Value = 310.0363185401719

$ws.Range("C95").Value = 30
$ws.Range("D95").Value = 265.7454158915759

$ws.Range("C96").Value = 31
$ws.Range("D96").Value = 274.6035964212951

# Rows 97-111 (old rows 96-110) already carry the correct shifted values
# and don't need any edits. Row 112 (old row 111, for 2021-02-27) now has
# enough data for the rolling window to be computable.
$ws.Range("C112").Value = 63
$ws.Range("D112").Value = 558.0653733723093

# One more day of data (2021-03-02) was appended at the bottom of the
# table; its rolling-window figures aren't computable yet so C/D stay
# blank, matching the existing trailing rows (113-114, old rows 112-113).
$ws.Range("A113").Copy()
$ws.Range("A115").PasteSpecial(-4122)
$ws.Range("A115").Value = 44257
$ws.Range("B115").Value = 5
